# Update market price / profit figures across the Leve profit worksheets.
# Generated from the scheduled market-data refresh diff.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @{Cell="H76"; Value=3079.14},
    @{Cell="I76"; Value=3012.0264},
    @{Cell="J76"; Value=3291.6667},
    @{Cell="K76"; Value=3012.0264},
    @{Cell="L76"; Value=3291.6667},
    @{Cell="M76"; Value=-2697.0264},
    @{Cell="N76"; Value=-3921.6667},
    @{Cell="H79"; Value=3079.14},
    @{Cell="I79"; Value=3012.0264},
    @{Cell="J79"; Value=3291.6667},
    @{Cell="K79"; Value=3012.0264},
    @{Cell="L79"; Value=3291.6667},
    @{Cell="M79"; Value=-1920.0264},
    @{Cell="N79"; Value=-5475.6667}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @{Cell="H32"; Value=10198.305},
    @{Cell="I32"; Value=10118.28},
    @{Cell="K32"; Value=10118.28},
    @{Cell="M32"; Value=-9831.280000000001},
    @{Cell="H61"; Value=6251102.5},
    @{Cell="I61"; Value=7813437.5},
    @{Cell="J61"; Value=1764.1875},
    @{Cell="K61"; Value=7813437.5},
    @{Cell="L61"; Value=1764.1875},
    @{Cell="M61"; Value=-7813225.5},
    @{Cell="N61"; Value=-2188.1875},
    @{Cell="H63"; Value=3167.2222},
    @{Cell="I63"; Value=3000.8333},
    @{Cell="K63"; Value=3000.8333},
    @{Cell="M63"; Value=-2314.8333},
    @{Cell="H66"; Value=3167.2222},
    @{Cell="I66"; Value=3000.8333},
    @{Cell="K66"; Value=15004.1665},
    @{Cell="M66"; Value=-11572.1665},
    @{Cell="H136"; Value=6251102.5},
    @{Cell="I136"; Value=7813437.5},
    @{Cell="J136"; Value=1764.1875},
    @{Cell="K136"; Value=23440312.5},
    @{Cell="L136"; Value=5292.5625},
    @{Cell="M136"; Value=-23437762.5},
    @{Cell="N136"; Value=-10392.5625}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @{Cell="H20"; Value=3229.75},
    @{Cell="I20"; Value=3219.75},
    @{Cell="J20"; Value=3249.75},
    @{Cell="K20"; Value=3219.75},
    @{Cell="L20"; Value=3249.75},
    @{Cell="M20"; Value=-2972.75},
    @{Cell="N20"; Value=-3743.75},
    @{Cell="H82"; Value=14146.667},
    @{Cell="J82"; Value=19583.334},
    @{Cell="L82"; Value=19583.334},
    @{Cell="N82"; Value=-20349.334},
    @{Cell="H85"; Value=14146.667},
    @{Cell="J85"; Value=19583.334},
    @{Cell="L85"; Value=19583.334},
    @{Cell="N85"; Value=-22235.334},
    @{Cell="H99"; Value=1153.6428},
    @{Cell="I99"; Value=1127},
    @{Cell="J99"; Value=1500},
    @{Cell="K99"; Value=1127},
    @{Cell="L99"; Value=1500},
    @{Cell="M99"; Value=371},
    @{Cell="N99"; Value=-4496}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @{Cell="H31"; Value=4331826.5},
    @{Cell="I31"; Value=2910.1833},
    @{Cell="J31"; Value=19610356},
    @{Cell="K31"; Value=2910.1833},
    @{Cell="L31"; Value=19610356},
    @{Cell="M31"; Value=-2615.1833},
    @{Cell="N31"; Value=-19610946},
    @{Cell="H34"; Value=4331826.5},
    @{Cell="I34"; Value=2910.1833},
    @{Cell="J34"; Value=19610356},
    @{Cell="K34"; Value=2910.1833},
    @{Cell="L34"; Value=19610356},
    @{Cell="M34"; Value=-2708.1833},
    @{Cell="N34"; Value=-19610760},
    @{Cell="H132"; Value=6251264.5},
    @{Cell="I132"; Value=7693293},
    @{Cell="J132"; Value=2474.1333},
    @{Cell="K132"; Value=23079879},
    @{Cell="L132"; Value=7422.3999},
    @{Cell="M132"; Value=-23077349},
    @{Cell="N132"; Value=-12482.3999},
    @{Cell="H134"; Value=1054.5344},
    @{Cell="I134"; Value=998.95746},
    @{Cell="J134"; Value=1292},
    @{Cell="K134"; Value=2996.87238},
    @{Cell="L134"; Value=3876},
    @{Cell="M134"; Value=-461.8723799999998},
    @{Cell="N134"; Value=-8946}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @{Cell="H5"; Value=3322.3333},
    @{Cell="I5"; Value=298.66666},
    @{Cell="J5"; Value=4834.1665},
    @{Cell="K5"; Value=895.9999799999999},
    @{Cell="L5"; Value=14502.4995},
    @{Cell="M5"; Value=-783.9999799999999},
    @{Cell="N5"; Value=-14726.4995},
    @{Cell="H106"; Value=10832.75},
    @{Cell="J106"; Value=10832.75},
    @{Cell="L106"; Value=32498.25},
    @{Cell="N106"; Value=-34390.25},
    @{Cell="H114"; Value=1638.3889},
    @{Cell="I114"; Value=1303},
    @{Cell="J114"; Value=2057.625},
    @{Cell="K114"; Value=3909},
    @{Cell="L114"; Value=6172.875},
    @{Cell="M114"; Value=-655},
    @{Cell="N114"; Value=-12680.875},
    @{Cell="H120"; Value=15322.167},
    @{Cell="I120"; Value=0},
    @{Cell="J120"; Value=15322.167},
    @{Cell="K120"; Value=0},
    @{Cell="L120"; Value=45966.501},
    @{Cell="N120"; Value=-55642.501},
    @{Cell="H135"; Value=3322.3333},
    @{Cell="I135"; Value=298.66666},
    @{Cell="J135"; Value=4834.1665},
    @{Cell="K135"; Value=2687.99994},
    @{Cell="L135"; Value=43507.4985},
    @{Cell="M135"; Value=-152.9999399999997},
    @{Cell="N135"; Value=-48577.4985}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
$ws.Range("M120").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @{Cell="H70"; Value=6759.9365},
    @{Cell="I70"; Value=8143.6743},
    @{Cell="J70"; Value=3784.9},
    @{Cell="K70"; Value=8143.6743},
    @{Cell="L70"; Value=3784.9},
    @{Cell="M70"; Value=-7873.6743},
    @{Cell="N70"; Value=-4324.9},
    @{Cell="H73"; Value=6759.9365},
    @{Cell="I73"; Value=8143.6743},
    @{Cell="J73"; Value=3784.9},
    @{Cell="K73"; Value=8143.6743},
    @{Cell="L73"; Value=3784.9},
    @{Cell="M73"; Value=-7207.6743},
    @{Cell="N73"; Value=-5656.9},
    @{Cell="H80"; Value=12349008},
    @{Cell="I80"; Value=33335570},
    @{Cell="J80"; Value=3970.9412},
    @{Cell="K80"; Value=33335570},
    @{Cell="L80"; Value=3970.9412},
    @{Cell="M80"; Value=-33334572},
    @{Cell="N80"; Value=-5966.9412},
    @{Cell="H83"; Value=12349008},
    @{Cell="I83"; Value=33335570},
    @{Cell="J83"; Value=3970.9412},
    @{Cell="K83"; Value=166677850},
    @{Cell="L83"; Value=19854.706},
    @{Cell="M83"; Value=-166672858},
    @{Cell="N83"; Value=-29838.706},
    @{Cell="H126"; Value=6342.857},
    @{Cell="I126"; Value=4000},
    @{Cell="J126"; Value=6733.3335},
    @{Cell="K126"; Value=12000},
    @{Cell="L126"; Value=20200.0005},
    @{Cell="M126"; Value=-9530},
    @{Cell="N126"; Value=-25140.0005},
    @{Cell="H132"; Value=3281.2407},
    @{Cell="I132"; Value=2340.7104},
    @{Cell="J132"; Value=5515},
    @{Cell="K132"; Value=7022.1312},
    @{Cell="L132"; Value=16545},
    @{Cell="M132"; Value=-4492.1312},
    @{Cell="N132"; Value=-21605}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @{Cell="H40"; Value=10374.875},
    @{Cell="I40"; Value=26000},
    @{Cell="J40"; Value=5166.5},
    @{Cell="K40"; Value=26000},
    @{Cell="L40"; Value=5166.5},
    @{Cell="M40"; Value=-25864},
    @{Cell="N40"; Value=-5438.5},
    @{Cell="H82"; Value=1918},
    @{Cell="I82"; Value=1875},
    @{Cell="J82"; Value=1942.5714},
    @{Cell="K82"; Value=1875},
    @{Cell="L82"; Value=1942.5714},
    @{Cell="M82"; Value=-1514},
    @{Cell="N82"; Value=-2664.5714},
    @{Cell="H85"; Value=1918},
    @{Cell="I85"; Value=1875},
    @{Cell="J85"; Value=1942.5714},
    @{Cell="K85"; Value=1875},
    @{Cell="L85"; Value=1942.5714},
    @{Cell="M85"; Value=-627},
    @{Cell="N85"; Value=-4438.5714},
    @{Cell="H136"; Value=12503767},
    @{Cell="I136"; Value=17242550},
    @{Cell="J136"; Value=10609.546},
    @{Cell="K136"; Value=51727650},
    @{Cell="L136"; Value=31828.638},
    @{Cell="M136"; Value=-51725100},
    @{Cell="N136"; Value=-36928.638}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @{Cell="H10"; Value=36756},
    @{Cell="I10"; Value=0},
    @{Cell="J10"; Value=36756},
    @{Cell="K10"; Value=0},
    @{Cell="L10"; Value=36756},
    @{Cell="N10"; Value=-37094},
    @{Cell="H13"; Value=2750},
    @{Cell="I13"; Value=500},
    @{Cell="J13"; Value=5000},
    @{Cell="K13"; Value=500},
    @{Cell="L13"; Value=5000},
    @{Cell="M13"; Value=-360},
    @{Cell="N13"; Value=-5280},
    @{Cell="H132"; Value=2713.7},
    @{Cell="I132"; Value=2058.875},
    @{Cell="J132"; Value=5333},
    @{Cell="K132"; Value=6176.625},
    @{Cell="L132"; Value=15999},
    @{Cell="M132"; Value=-3646.625},
    @{Cell="N132"; Value=-21059},
    @{Cell="H136"; Value=1082.2122},
    @{Cell="I136"; Value=903.7143},
    @{Cell="J136"; Value=2081.8},
    @{Cell="K136"; Value=2711.1429},
    @{Cell="L136"; Value=6245.400000000001},
    @{Cell="M136"; Value=-161.1428999999998},
    @{Cell="N136"; Value=-11345.4}
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
$ws.Range("M10").ClearContents()

Write-Host "Applied Ultima_Profits market data refresh."